$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update sheet (tab) name and title to reflect new "through" date
$ws.Name = "Through 2022-03-11"

# Update the "March (through ...)" label in column A, row 4
$ws.Range("A4").Value = "March (through 03-11)"

# Row 3 - February: only 2022 (I) column changes
$ws.Range("I3").Value = 141

# Row 4 - March: update all year columns (B..I)
$ws.Range("B4").Value = 11
$ws.Range("C4").Value = 15
$ws.Range("D4").Value = 20
$ws.Range("E4").Value = 20
$ws.Range("F4").Value = 12
$ws.Range("G4").Value = 20
$ws.Range("H4").Value = 32
$ws.Range("I4").Value = 48

# Row 5 - Total: update all year columns (B..I)
$ws.Range("B5").Value = 48
$ws.Range("C5").Value = 102
$ws.Range("D5").Value = 151
$ws.Range("E5").Value = 157
$ws.Range("F5").Value = 91
$ws.Range("G5").Value = 161
$ws.Range("H5").Value = 374
$ws.Range("I5").Value = 348
